# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit figures across the per-class leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1127.6945
$ws.Range("I15").Value = 1127.6945
$ws.Range("K15").Value = 3383.0835
$ws.Range("M15").Value = -3214.0835

$ws.Range("H28").Value = 3298.9092
$ws.Range("I28").Value = 367.14285
$ws.Range("K28").Value = 367.14285
$ws.Range("M28").Value = 117.85715

$ws.Range("H42").Value = 933.5333000000001
$ws.Range("J42").Value = 1621.25
$ws.Range("L42").Value = 4863.75
$ws.Range("N42").Value = -5323.75

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()

$ws.Range("H69").Value = 13582.5
$ws.Range("I69").Value = 6666
$ws.Range("K69").Value = 19998
$ws.Range("M69").Value = -19124

$ws.Range("H72").Value = 13582.5
$ws.Range("I72").Value = 6666
$ws.Range("K72").Value = 59994
$ws.Range("M72").Value = -55626

$ws.Range("H82").Value = 7820.4
$ws.Range("I82").Value = 4707.6665
$ws.Range("J82").Value = 12489.5
$ws.Range("K82").Value = 14122.9995
$ws.Range("L82").Value = 37468.5
$ws.Range("M82").Value = -13716.9995
$ws.Range("N82").Value = -38280.5

$ws.Range("H85").Value = 7820.4
$ws.Range("I85").Value = 4707.6665
$ws.Range("J85").Value = 12489.5
$ws.Range("K85").Value = 14122.9995
$ws.Range("L85").Value = 37468.5
$ws.Range("M85").Value = -12718.9995
$ws.Range("N85").Value = -40276.5

$ws.Range("H92").Value = 871.1667
$ws.Range("I92").Value = 393.625
$ws.Range("J92").Value = 1826.25
$ws.Range("K92").Value = 393.625
$ws.Range("L92").Value = 1826.25
$ws.Range("M92").Value = 854.375
$ws.Range("N92").Value = -4322.25

$ws.Range("H100").Value = 11354.182
$ws.Range("I100").Value = 7159.6
$ws.Range("J100").Value = 14849.667
$ws.Range("K100").Value = 7159.6
$ws.Range("L100").Value = 14849.667
$ws.Range("M100").Value = -6618.6
$ws.Range("N100").Value = -15931.667

$ws.Range("H101").Value = 2756.6
$ws.Range("J101").Value = 3867.5715
$ws.Range("L101").Value = 11602.7145
$ws.Range("N101").Value = -14846.7145

$ws.Range("H131").Value = 1644217.2
$ws.Range("J131").Value = 1770388.1
$ws.Range("L131").Value = 5311164.300000001
$ws.Range("N131").Value = -5321244.300000001

$ws.Range("H135").Value = 3238.2942
$ws.Range("I135").Value = 1383.3334
$ws.Range("J135").Value = 5325.125
$ws.Range("K135").Value = 12450.0006
$ws.Range("L135").Value = 47926.125
$ws.Range("M135").Value = -9915.000599999999
$ws.Range("N135").Value = -52996.125

$ws.Range("H138").Value = 3215
$ws.Range("J138").Value = 2189.5
$ws.Range("L138").Value = 6568.5
$ws.Range("N138").Value = -16848.5

$ws.Range("H141").Value = 15642943
$ws.Range("I141").Value = 19235338
$ws.Range("J141").Value = 75896
$ws.Range("K141").Value = 57706014
$ws.Range("L141").Value = 227688
$ws.Range("M141").Value = -57700834
$ws.Range("N141").Value = -238048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5128.769
$ws.Range("I110").Value = 6220.6665
$ws.Range("K110").Value = 6220.6665
$ws.Range("M110").Value = -4175.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1992.88
$ws.Range("I94").Value = 2100.6667
$ws.Range("K94").Value = 2100.6667
$ws.Range("M94").Value = -1649.6667

$ws.Range("H99").Value = 2851.4285
$ws.Range("J99").Value = 2400
$ws.Range("L99").Value = 2400
$ws.Range("N99").Value = -5396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2718.9092
$ws.Range("I58").Value = 2760.4285
$ws.Range("K58").Value = 2760.4285
$ws.Range("M58").Value = -2557.4285

$ws.Range("H122").Value = 5146.9165
$ws.Range("I122").Value = 5470.5
$ws.Range("K122").Value = 16411.5
$ws.Range("M122").Value = -13961.5

$ws.Range("H132").Value = 3063.9524
$ws.Range("I132").Value = 2868.75
$ws.Range("J132").Value = 3324.2222
$ws.Range("K132").Value = 8606.25
$ws.Range("L132").Value = 9972.6666
$ws.Range("M132").Value = -6076.25
$ws.Range("N132").Value = -15032.6666

$ws.Range("H136").Value = 2718.9092
$ws.Range("I136").Value = 2760.4285
$ws.Range("K136").Value = 8281.2855
$ws.Range("M136").Value = -5731.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1949
$ws.Range("I5").Value = 1949
$ws.Range("K5").Value = 5847
$ws.Range("M5").Value = -5735

$ws.Range("H68").Value = 2943.875
$ws.Range("I68").Value = 2650.2856
$ws.Range("J68").Value = 4999
$ws.Range("K68").Value = 7950.8568
$ws.Range("L68").Value = 14997
$ws.Range("M68").Value = -7139.8568
$ws.Range("N68").Value = -16619

$ws.Range("H71").Value = 2943.875
$ws.Range("I71").Value = 2650.2856
$ws.Range("J71").Value = 4999
$ws.Range("K71").Value = 23852.5704
$ws.Range("L71").Value = 44991
$ws.Range("M71").Value = -19796.5704
$ws.Range("N71").Value = -53103

$ws.Range("H135").Value = 1949
$ws.Range("I135").Value = 1949
$ws.Range("K135").Value = 17541
$ws.Range("M135").Value = -15006

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4169113.8
$ws.Range("I68").Value = 10418452
$ws.Range("J68").Value = 2888.5
$ws.Range("K68").Value = 10418452
$ws.Range("L68").Value = 2888.5
$ws.Range("M68").Value = -10417703
$ws.Range("N68").Value = -4386.5

$ws.Range("H71").Value = 4169113.8
$ws.Range("I71").Value = 10418452
$ws.Range("J71").Value = 2888.5
$ws.Range("K71").Value = 52092260
$ws.Range("L71").Value = 14442.5
$ws.Range("M71").Value = -52088516
$ws.Range("N71").Value = -21930.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 687.625
$ws.Range("I113").Value = 628.5
$ws.Range("K113").Value = 1885.5
$ws.Range("M113").Value = 284.5

$ws.Range("H132").Value = 594485.4399999999
$ws.Range("I132").Value = 5768.909
$ws.Range("J132").Value = 1673799
$ws.Range("K132").Value = 17306.727
$ws.Range("L132").Value = 5021397
$ws.Range("M132").Value = -14776.727
$ws.Range("N132").Value = -5026457
